$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MACRO_SCORE column (N) for rows 2-7 with the refreshed value.
$ws.Range("N2:N7").Value = 85.92117485762657
